$wb = $excel.ActiveWorkbook

# --- Remove header-row styling (bold font + thin border + centered alignment) ---
# Revert A1:N1 on every sheet back to the default 'Normal' style.
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1:N1").Style = "Normal"
}

# --- Update recalculated market-price figures per sheet ---

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1694.4762
$ws.Range("J19").Value = 1779.75
$ws.Range("L19").Value = 1779.75
$ws.Range("N19").Value = -2129.75
$ws.Range("H43").Value = 1960
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1960
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1960
$ws.Range("M43").ClearContents()  # was -931
$ws.Range("N43").Value = -2098
$ws.Range("H64").Value = 2319.1667
$ws.Range("I64").Value = 2088.889
$ws.Range("K64").Value = 2088.889
$ws.Range("M64").Value = -1840.889
$ws.Range("H67").Value = 2319.1667
$ws.Range("I67").Value = 2088.889
$ws.Range("K67").Value = 2088.889
$ws.Range("M67").Value = -1230.889
$ws.Range("H113").Value = 2976.2104
$ws.Range("I113").Value = 2972.625
$ws.Range("K113").Value = 2972.625
$ws.Range("M113").Value = 281.375
$ws.Range("H127").Value = 2015.3
$ws.Range("I127").Value = 2156.625
$ws.Range("J127").Value = 1450
$ws.Range("K127").Value = 6469.875
$ws.Range("L127").Value = 4350
$ws.Range("M127").Value = -1509.875
$ws.Range("N127").Value = -14270
$ws.Range("H131").Value = 2474.9333
$ws.Range("J131").Value = 2740.6924
$ws.Range("L131").Value = 8222.0772
$ws.Range("N131").Value = -18302.0772
$ws.Range("H132").Value = 1328.9231
$ws.Range("I132").Value = 1328.9231
$ws.Range("K132").Value = 3986.7693
$ws.Range("M132").Value = -1456.7693
$ws.Range("H138").Value = 3503.2
$ws.Range("I138").Value = 4431.5557
$ws.Range("J138").Value = 2110.6667
$ws.Range("K138").Value = 13294.6671
$ws.Range("L138").Value = 6332.000100000001
$ws.Range("M138").Value = -8154.667099999999
$ws.Range("N138").Value = -16612.0001
$ws.Range("H141").Value = 967884.25
$ws.Range("I141").Value = 1335268.9
$ws.Range("J141").Value = 3499.625
$ws.Range("K141").Value = 4005806.7
$ws.Range("L141").Value = 10498.875
$ws.Range("M141").Value = -4000626.7
$ws.Range("N141").Value = -20858.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2315.2917
$ws.Range("I32").Value = 1598.322
$ws.Range("K32").Value = 1598.322
$ws.Range("M32").Value = -1311.322
$ws.Range("H45").Value = 2648.9119
$ws.Range("I45").Value = 2484.3333
$ws.Range("K45").Value = 2484.3333
$ws.Range("M45").Value = -2107.3333
$ws.Range("H61").Value = 3426.5
$ws.Range("I61").Value = 2651.8462
$ws.Range("K61").Value = 2651.8462
$ws.Range("M61").Value = -2439.8462
$ws.Range("H74").Value = 1693.2142
$ws.Range("I74").Value = 891.36365
$ws.Range("J74").Value = 4633.3335
$ws.Range("K74").Value = 891.36365
$ws.Range("L74").Value = 4633.3335
$ws.Range("M74").Value = -17.36365000000001
$ws.Range("N74").Value = -6381.3335
$ws.Range("H77").Value = 1693.2142
$ws.Range("I77").Value = 891.36365
$ws.Range("J77").Value = 4633.3335
$ws.Range("K77").Value = 4456.81825
$ws.Range("L77").Value = 23166.6675
$ws.Range("M77").Value = -88.81825000000026
$ws.Range("N77").Value = -31902.6675
$ws.Range("H122").Value = 127530.664
$ws.Range("H132").Value = 2579.7727
$ws.Range("I132").Value = 2371.3713
$ws.Range("K132").Value = 7114.113899999999
$ws.Range("M132").Value = -4584.113899999999
$ws.Range("H136").Value = 3426.5
$ws.Range("I136").Value = 2651.8462
$ws.Range("K136").Value = 7955.5386
$ws.Range("M136").Value = -5405.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1360.6451
$ws.Range("I94").Value = 1163.5
$ws.Range("K94").Value = 1163.5
$ws.Range("M94").Value = -712.5
$ws.Range("H105").Value = 2397.625
$ws.Range("I105").Value = 2401.9546
$ws.Range("K105").Value = 2401.9546
$ws.Range("M105").Value = -654.9546
$ws.Range("H134").Value = 3038.4
$ws.Range("I134").Value = 2814.75
$ws.Range("J134").Value = 3933
$ws.Range("K134").Value = 8444.25
$ws.Range("L134").Value = 11799
$ws.Range("M134").Value = -5909.25
$ws.Range("N134").Value = -16869

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7813073.5
$ws.Range("I22").Value = 570.5
$ws.Range("J22").Value = 15625576
$ws.Range("K22").Value = 570.5
$ws.Range("L22").Value = 15625576
$ws.Range("M22").Value = -220.5
$ws.Range("N22").Value = -15626276
$ws.Range("H31").Value = 1547.2
$ws.Range("I31").Value = 902.65
$ws.Range("J31").Value = 2062.84
$ws.Range("K31").Value = 902.65
$ws.Range("L31").Value = 2062.84
$ws.Range("M31").Value = -607.65
$ws.Range("N31").Value = -2652.84
$ws.Range("H34").Value = 1547.2
$ws.Range("I34").Value = 902.65
$ws.Range("J34").Value = 2062.84
$ws.Range("K34").Value = 902.65
$ws.Range("L34").Value = 2062.84
$ws.Range("M34").Value = -700.65
$ws.Range("N34").Value = -2466.84
$ws.Range("H58").Value = 2718738
$ws.Range("I58").Value = 3953551.8
$ws.Range("J58").Value = 2148
$ws.Range("K58").Value = 3953551.8
$ws.Range("L58").Value = 2148
$ws.Range("M58").Value = -3953348.8
$ws.Range("N58").Value = -2554
$ws.Range("H94").Value = 770.1
$ws.Range("I94").Value = 668.7778
$ws.Range("J94").Value = 853
$ws.Range("K94").Value = 668.7778
$ws.Range("L94").Value = 853
$ws.Range("M94").Value = -217.7778
$ws.Range("N94").Value = -1755
$ws.Range("H99").Value = 2999.8572
$ws.Range("I99").Value = 2499.6667
$ws.Range("K99").Value = 2499.6667
$ws.Range("M99").Value = -1001.6667
$ws.Range("H122").Value = 1022.6
$ws.Range("I122").Value = 1105.5454
$ws.Range("K122").Value = 3316.6362
$ws.Range("M122").Value = -866.6361999999999
$ws.Range("H126").Value = 2999.8572
$ws.Range("I126").Value = 2499.6667
$ws.Range("K126").Value = 7499.000100000001
$ws.Range("M126").Value = -5029.000100000001
$ws.Range("H132").Value = 2815.75
$ws.Range("I132").Value = 2308.2632
$ws.Range("J132").Value = 3887.111
$ws.Range("K132").Value = 6924.7896
$ws.Range("L132").Value = 11661.333
$ws.Range("M132").Value = -4394.7896
$ws.Range("N132").Value = -16721.333
$ws.Range("H136").Value = 2718738
$ws.Range("I136").Value = 3953551.8
$ws.Range("J136").Value = 2148
$ws.Range("K136").Value = 11860655.4
$ws.Range("L136").Value = 6444
$ws.Range("M136").Value = -11858105.4
$ws.Range("N136").Value = -11544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 23810110
$ws.Range("I114").Value = 705.4
$ws.Range("J114").Value = 142857140
$ws.Range("K114").Value = 2116.2
$ws.Range("L114").Value = 428571420
$ws.Range("M114").Value = 1137.8
$ws.Range("N114").Value = -428577928
$ws.Range("H121").Value = 681.2222
$ws.Range("J121").Value = 747.4286
$ws.Range("L121").Value = 2242.2858
$ws.Range("N121").Value = -4862.2858
$ws.Range("H129").Value = 52530.285
$ws.Range("I129").Value = 692.8333
$ws.Range("J129").Value = 91408.375
$ws.Range("K129").Value = 2078.4999
$ws.Range("L129").Value = 274225.125
$ws.Range("M129").Value = 2921.5001
$ws.Range("N129").Value = -284225.125
$ws.Range("H131").Value = 9104736
$ws.Range("J131").Value = 14329.528
$ws.Range("L131").Value = 42988.584
$ws.Range("N131").Value = -53068.584
$ws.Range("H140").Value = 2888.1667
$ws.Range("I140").Value = 866.8
$ws.Range("J140").Value = 4332
$ws.Range("K140").Value = 2600.4
$ws.Range("L140").Value = 12996
$ws.Range("M140").Value = 2579.6
$ws.Range("N140").Value = -23356

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1009.7
$ws.Range("I97").Value = 771.35
$ws.Range("J97").Value = 1486.4
$ws.Range("K97").Value = 771.35
$ws.Range("L97").Value = 1486.4
$ws.Range("M97").Value = -275.35
$ws.Range("N97").Value = -2478.4
$ws.Range("H102").Value = 3052.5454
$ws.Range("I102").Value = 3158.4
$ws.Range("K102").Value = 3158.4
$ws.Range("M102").Value = -1536.4
$ws.Range("H122").Value = 3257.125
$ws.Range("I122").Value = 1366
$ws.Range("K122").Value = 4098
$ws.Range("M122").Value = -1648

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4059.842
$ws.Range("I7").Value = 3041.6
$ws.Range("J7").Value = 4423.5
$ws.Range("K7").Value = 3041.6
$ws.Range("L7").Value = 4423.5
$ws.Range("M7").Value = -2929.6
$ws.Range("N7").Value = -4647.5
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10348
$ws.Range("H22").Value = 4095.7
$ws.Range("J22").Value = 4494.625
$ws.Range("L22").Value = 4494.625
$ws.Range("N22").Value = -5084.625
$ws.Range("H27").Value = 4095.7
$ws.Range("J27").Value = 4494.625
$ws.Range("L27").Value = 4494.625
$ws.Range("N27").Value = -4708.625
$ws.Range("H40").Value = 7099
$ws.Range("I40").Value = 3148.875
$ws.Range("J40").Value = 14999.25
$ws.Range("K40").Value = 3148.875
$ws.Range("L40").Value = 14999.25
$ws.Range("M40").Value = -3012.875
$ws.Range("N40").Value = -15271.25
$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1400
$ws.Range("K46").Value = 1400
$ws.Range("M46").Value = -1212
$ws.Range("H110").Value = 25000
$ws.Range("J110").Value = 25000
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180
$ws.Range("H122").Value = 8032.2856
$ws.Range("I122").Value = 5633.9287
$ws.Range("J122").Value = 12829
$ws.Range("K122").Value = 16901.7861
$ws.Range("L122").Value = 38487
$ws.Range("M122").Value = -14451.7861
$ws.Range("N122").Value = -43387
$ws.Range("H126").Value = 4059.842
$ws.Range("I126").Value = 3041.6
$ws.Range("J126").Value = 4423.5
$ws.Range("K126").Value = 9124.799999999999
$ws.Range("L126").Value = 13270.5
$ws.Range("M126").Value = -6654.799999999999
$ws.Range("N126").Value = -18210.5
$ws.Range("H132").Value = 3311.6538
$ws.Range("I132").Value = 1151
$ws.Range("K132").Value = 3453
$ws.Range("M132").Value = -923
$ws.Range("H136").Value = 4869.5
$ws.Range("I136").Value = 2905.6365
$ws.Range("J136").Value = 7955.5713
$ws.Range("K136").Value = 8716.9095
$ws.Range("L136").Value = 23866.7139
$ws.Range("M136").Value = -6166.9095
$ws.Range("N136").Value = -28966.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45397.586
$ws.Range("I122").Value = 48653
$ws.Range("K122").Value = 145959
$ws.Range("M122").Value = -143509
$ws.Range("H123").Value = 39090.625
$ws.Range("J123").Value = 39090.625
$ws.Range("L123").Value = 39090.625
$ws.Range("N123").Value = -48890.625
$ws.Range("H132").Value = 3614.077
$ws.Range("I132").Value = 3117.6
$ws.Range("J132").Value = 3924.375
$ws.Range("K132").Value = 9352.799999999999
$ws.Range("L132").Value = 11773.125
$ws.Range("M132").Value = -6822.799999999999
$ws.Range("N132").Value = -16833.125
$ws.Range("H136").Value = 24157166
$ws.Range("I136").Value = 30866466
$ws.Range("K136").Value = 92599398
$ws.Range("M136").Value = -92596848
